$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F3").Value = 1242
    $ws.Range("F5").Value = 1
    $ws.Range("F7").Value = 163
}
